$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige el nombre mal escrito "Corsea del Sur" -> "Corea del Sur"
$ws.Range("B7").Value = "Corea del Sur"
$ws.Range("C8").Value = "Corea del Sur"

# Ajusta la celda seleccionada
$ws.Range("B8").Select()
